$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" ---
# Overview sheet: columns E (zh-cn) and F (de-de), row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column (C), row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: Status column (C), row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Shrink the now-narrower status columns ---
# ColumnWidth is expressed in characters; the stored width snaps to the
# nearest whole pixel on the Calibri-11 grid (pixels = round(width*6)+5),
# so 12.576851254417766 is the input that lands on the target ~13.41 width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.576851254417766
$wsOverview.Columns.Item(6).ColumnWidth = 12.576851254417766
$wsZhCn.Columns.Item(3).ColumnWidth = 12.576851254417766
$wsDeDe.Columns.Item(3).ColumnWidth = 12.576851254417766
